{"js": "// Add a new \"Topic 3\" block (heading + File line + divider + blank spacer\n// lines) right after the existing \"Topic 2\" block, mirroring the layout\n// already used for \"Topic 1\" / \"Topic 2\".\n//\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Topic 2 : Working with Java Primitive types.\" heading\n// paragraph so the insertion point is found robustly (not by a hard-coded\n// index) even if earlier content shifts.\nlet topic2HeadingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Topic 2\") === 0 && t.indexOf(\"Working with Java Primitive types.\") !== -1) {\n    topic2HeadingIndex = i;\n    break;\n  }\n}\nif (topic2HeadingIndex === -1) {\n  throw new Error(\"Could not find the 'Topic 2' heading paragraph.\");\n}\n\n// Layout for the Topic 2 block is: heading, File line, divider. The new\n// \"Topic 3\" block is inserted right after that divider paragraph.\nconst dividerIndex = topic2HeadingIndex + 2;\nlet anchor = paragraphs.items[dividerIndex];\nanchor.load(\"text\");\nawait context.sync();\nif (anchor.text.indexOf(\"---\") !== 0) {\n  throw new Error(\"Unexpected paragraph at the 'Topic 2' divider position.\");\n}\n\n// Bold, 14pt (w:sz 28 half-points) run/paragraph-mark formatting shared by\n// every paragraph in this block.\nconst rPr = '<w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr>';\nconst pPr = `<w:pPr>${rPr}</w:pPr>`;\n\nfunction packageOoxml(fragmentXml) {\n  // Wrap a <w:p> fragment in the minimal flat-OPC \"pkg:package\" envelope\n  // insertOoxml() expects.\n  return '<?xml version=\"1.0\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    `<w:body>${fragmentXml}<w:sectPr/></w:body>` +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n}\n\nconst dashes = \"-\".repeat(105);\n\n// The 9 new paragraphs, in document order.\nconst newParagraphFragments = [\n  // 1. blank spacer line\n  `<w:p>${pPr}</w:p>`,\n\n  // 2. \"Topic 3 \\t: Working with Java Arrays.\" (title underlined)\n  `<w:p>${pPr}` +\n    `<w:r>${rPr}<w:t xml:space=\"preserve\">Topic 3 </w:t></w:r>` +\n    `<w:r>${rPr}<w:tab/><w:t xml:space=\"preserve\">: </w:t></w:r>` +\n    '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:u w:val=\"single\"/></w:rPr>' +\n    '<w:t>Working with Java Arrays.</w:t></w:r>' +\n  `</w:p>`,\n\n  // 3. \"File \\t\\t: C:\\Users\\Gaurav\\Desktop\\OCAJP SE 11 Certification/Topic 3.txt\"\n  `<w:p>${pPr}` +\n    `<w:r>${rPr}<w:t xml:space=\"preserve\">File </w:t></w:r>` +\n    `<w:r>${rPr}<w:tab/></w:r>` +\n    `<w:r>${rPr}<w:tab/><w:t xml:space=\"preserve\">: </w:t></w:r>` +\n    `<w:r>${rPr}<w:t xml:space=\"preserve\">C:\\\\Users\\\\Gaurav\\\\Desktop\\\\OCAJP SE 11 Certification/Topic </w:t></w:r>` +\n    `<w:r>${rPr}<w:t>3</w:t></w:r>` +\n    `<w:r>${rPr}<w:t>.txt</w:t></w:r>` +\n  `</w:p>`,\n\n  // 4. divider\n  `<w:p>${pPr}<w:r>${rPr}<w:t>${dashes}</w:t></w:r></w:p>`,\n\n  // 5-9. five blank spacer lines\n  `<w:p>${pPr}</w:p>`,\n  `<w:p>${pPr}</w:p>`,\n  `<w:p>${pPr}</w:p>`,\n  `<w:p>${pPr}</w:p>`,\n  `<w:p>${pPr}</w:p>`,\n];\n\nfor (const fragment of newParagraphFragments) {\n  // Create an (empty) paragraph right after the current anchor so a real\n  // paragraph boundary exists, then replace its contents with the exact\n  // OOXML fragment (this preserves literal <w:tab/> elements, which a\n  // plain insertText(\"\\t\") call would otherwise collapse into the text\n  // run as a raw tab character).\n  const newParagraph = anchor.insertParagraph(\"\", Word.InsertLocation.after);\n  await context.sync();\n\n  newParagraph.insertOoxml(packageOoxml(fragment), Word.InsertLocation.replace);\n  await context.sync();\n\n  anchor = newParagraph;\n}\n", "ps1": "# Add a new \"Topic 3\" block (heading + File line + divider + blank spacer\n# lines) right after the existing \"Topic 2\" block, mirroring the layout\n# already used for \"Topic 1\" / \"Topic 2\".\n\n$d = $word.ActiveDocument\n\n# Locate the \"Topic 2 : Working with Java Primitive types.\" heading\n# paragraph robustly (not by a hard-coded index).\n$topic2HeadingIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($t -like \"Topic 2*Working with Java Primitive types.*\") {\n        $topic2HeadingIndex = $i\n        break\n    }\n}\nif ($topic2HeadingIndex -eq -1) {\n    throw \"Could not find the 'Topic 2' heading paragraph.\"\n}\n\n# Layout for the Topic 2 block is: heading, File line, divider. The new\n# \"Topic 3\" block is inserted right after that divider paragraph.\n$dividerIndex = $topic2HeadingIndex + 2\n$anchorRange = $d.Paragraphs($dividerIndex).Range\nif ($anchorRange.Text.Substring(0, 3) -ne \"---\") {\n    throw \"Unexpected paragraph at the 'Topic 2' divider position.\"\n}\n\n# Bold, 14pt (w:sz 28 half-points) run/paragraph-mark formatting shared by\n# every paragraph in this block.\n$rPr = '<w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr>'\n$pPr = \"<w:pPr>$rPr</w:pPr>\"\n\nfunction New-PackageOoxml([string]$fragmentXml) {\n    # Wrap a <w:p> fragment in the minimal flat-OPC \"pkg:package\" envelope\n    # Range.InsertXML() expects.\n    return '<?xml version=\"1.0\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        \"<w:body>$fragmentXml<w:sectPr/></w:body>\" +\n        '</w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n$dashes = '---------------------------------------------------------------------------------------------------------'\n\n# The 9 new paragraphs, in document order.\n$newParagraphFragments = @(\n    # 1. blank spacer line\n    \"<w:p>$pPr</w:p>\",\n\n    # 2. \"Topic 3 `t: Working with Java Arrays.\" (title underlined)\n    (\"<w:p>$pPr\" +\n        \"<w:r>$rPr<w:t xml:space=`\"preserve`\">Topic 3 </w:t></w:r>\" +\n        \"<w:r>$rPr<w:tab/><w:t xml:space=`\"preserve`\">: </w:t></w:r>\" +\n        '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:u w:val=\"single\"/></w:rPr>' +\n        '<w:t>Working with Java Arrays.</w:t></w:r>' +\n        \"</w:p>\"),\n\n    # 3. \"File `t`t: C:\\Users\\Gaurav\\Desktop\\OCAJP SE 11 Certification/Topic 3.txt\"\n    (\"<w:p>$pPr\" +\n        \"<w:r>$rPr<w:t xml:space=`\"preserve`\">File </w:t></w:r>\" +\n        \"<w:r>$rPr<w:tab/></w:r>\" +\n        \"<w:r>$rPr<w:tab/><w:t xml:space=`\"preserve`\">: </w:t></w:r>\" +\n        \"<w:r>$rPr<w:t xml:space=`\"preserve`\">C:\\Users\\Gaurav\\Desktop\\OCAJP SE 11 Certification/Topic </w:t></w:r>\" +\n        \"<w:r>$rPr<w:t>3</w:t></w:r>\" +\n        \"<w:r>$rPr<w:t>.txt</w:t></w:r>\" +\n        \"</w:p>\"),\n\n    # 4. divider\n    \"<w:p>$pPr<w:r>$rPr<w:t>$dashes</w:t></w:r></w:p>\",\n\n    # 5-9. five blank spacer lines\n    \"<w:p>$pPr</w:p>\",\n    \"<w:p>$pPr</w:p>\",\n    \"<w:p>$pPr</w:p>\",\n    \"<w:p>$pPr</w:p>\",\n    \"<w:p>$pPr</w:p>\"\n)\n\n$anchorIndex = $dividerIndex\nforeach ($fragment in $newParagraphFragments) {\n    # Create an (empty) paragraph right after the current anchor so a real\n    # paragraph boundary exists, then inject the exact OOXML fragment into\n    # it (this preserves literal <w:tab/> elements, which plain text\n    # insertion of a tab character would otherwise collapse into the run\n    # as a raw tab character instead of a dedicated element).\n    $anchorRange = $d.Paragraphs($anchorIndex).Range\n    $anchorRange.InsertParagraphAfter()\n\n    $newIndex = $anchorIndex + 1\n    $newRange = $d.Paragraphs($newIndex).Range\n    $newRange.InsertXML((New-PackageOoxml $fragment))\n\n    $anchorIndex = $newIndex\n}\n"}
